$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 111936866
$ws.Range("B10").Value = 89557
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 5432
$ws.Range("F10").Value = "Granticka"
$ws.Range("G10").Value = "Porodaedalea chrysoloma"
$ws.Range("H10").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q10").Value = 448766
$ws.Range("R10").Value = 7087417
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("AC10").ClearContents()

$ws.Range("A11").Value = 111936792
$ws.Range("B11").Value = 90221
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 3298
$ws.Range("F11").Value = "Trådticka"
$ws.Range("G11").Value = "Climacocystis borealis"
$ws.Range("H11").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q11").Value = 448761
$ws.Range("R11").Value = 7087579
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("AC11").ClearContents()

$ws.Range("A12").Value = 111936865
$ws.Range("B12").Value = 89557
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 5432
$ws.Range("F12").Value = "Granticka"
$ws.Range("G12").Value = "Porodaedalea chrysoloma"
$ws.Range("H12").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q12").Value = 448738
$ws.Range("R12").Value = 7087426
$ws.Range("K12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("AC12").ClearContents()

$ws.Range("A13").Value = 111936893
$ws.Range("B13").Value = 77636
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("Q13").Value = 448742
$ws.Range("R13").Value = 7087502
$ws.Range("K13").ClearContents()
$ws.Range("L13").ClearContents()
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("AC13").ClearContents()

$ws.Range("A14").Value = 111936796
$ws.Range("B14").Value = 56430
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 100109
$ws.Range("F14").Value = "Tretåig hackspett"
$ws.Range("G14").Value = "Picoides tridactylus"
$ws.Range("H14").Value = "(Linnaeus, 1758)"
$ws.Range("Q14").Value = 448883
$ws.Range("R14").Value = 7087229
$ws.Range("K14").Value = ""
$ws.Range("L14").Value = ""
$ws.Range("M14").Value = ""
$ws.Range("N14").Value = ""
$ws.Range("AC14").Value = "ringhack äldre"

$ws.Range("A15").Value = 111936798
$ws.Range("B15").Value = 56430
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 100109
$ws.Range("F15").Value = "Tretåig hackspett"
$ws.Range("G15").Value = "Picoides tridactylus"
$ws.Range("H15").Value = "(Linnaeus, 1758)"
$ws.Range("Q15").Value = 448923
$ws.Range("R15").Value = 7087371
$ws.Range("K15").Value = ""
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = ""
$ws.Range("N15").Value = ""
$ws.Range("AC15").Value = "ringhack äldre"

$ws.Range("A16").Value = 111936870
$ws.Range("B16").Value = 89557
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 5432
$ws.Range("F16").Value = "Granticka"
$ws.Range("G16").Value = "Porodaedalea chrysoloma"
$ws.Range("H16").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q16").Value = 449019
$ws.Range("R16").Value = 7087277
$ws.Range("K16").ClearContents()
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("AC16").ClearContents()

$ws.Range("A17").Value = 111936795
$ws.Range("B17").Value = 56430
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 100109
$ws.Range("F17").Value = "Tretåig hackspett"
$ws.Range("G17").Value = "Picoides tridactylus"
$ws.Range("H17").Value = "(Linnaeus, 1758)"
$ws.Range("Q17").Value = 448749
$ws.Range("R17").Value = 7087422
$ws.Range("K17").Value = ""
$ws.Range("L17").Value = ""
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = ""
$ws.Range("AC17").Value = "ringhack äldre"

$ws.Range("A18").Value = 111936858
$ws.Range("B18").Value = 89979
$ws.Range("D18").Value = "VU"
$ws.Range("E18").Value = 1209
$ws.Range("F18").Value = "Rynkskinn"
$ws.Range("G18").Value = "Phlebia centrifuga"
$ws.Range("H18").Value = "P.Karst."
$ws.Range("Q18").Value = 448737
$ws.Range("R18").Value = 7087496
$ws.Range("K18").ClearContents()
$ws.Range("L18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("AC18").ClearContents()

$ws.Range("A19").Value = 111936869
$ws.Range("B19").Value = 89557
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 5432
$ws.Range("F19").Value = "Granticka"
$ws.Range("G19").Value = "Porodaedalea chrysoloma"
$ws.Range("H19").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q19").Value = 449144
$ws.Range("R19").Value = 7087118
$ws.Range("K19").ClearContents()
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("AC19").ClearContents()

$ws.Range("A20").Value = 111936867
$ws.Range("B20").Value = 89557
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 5432
$ws.Range("F20").Value = "Granticka"
$ws.Range("G20").Value = "Porodaedalea chrysoloma"
$ws.Range("H20").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q20").Value = 448792
$ws.Range("R20").Value = 7087386
$ws.Range("K20").ClearContents()
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("AC20").ClearContents()

$ws.Range("A21").Value = 111936868
$ws.Range("B21").Value = 89557
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 5432
$ws.Range("F21").Value = "Granticka"
$ws.Range("G21").Value = "Porodaedalea chrysoloma"
$ws.Range("H21").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q21").Value = 448988
$ws.Range("R21").Value = 7087187
$ws.Range("K21").ClearContents()
$ws.Range("L21").ClearContents()
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("AC21").ClearContents()
